$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: "1020. Number of Enclaves" (graph theory problem)
$row = 13

$ws.Range("A$row").Value = "1020. Number of Enclaves"

# Difficulty column re-uses the same look (fill) as the rest of column B ("Medium")
$ws.Range("B12").Copy($ws.Range("B$row"))
$ws.Range("B$row").Value = "Medium"

$ws.Range("C$row").Value = "Matrix Graphs"

$ws.Range("D$row").Value = "Flood fill from the edges. DFS from any 1s on the edges and mark them as 0, then to another pass to count the remaining 1s."

$url = "https://leetcode.com/problems/number-of-enclaves/solutions/3388131/python-java-c-simple-solution-easy-to-understand/?envType=study-plan-v2&envId=graph-theory "
$ws.Range("E$row").Value = $url
$ws.Hyperlinks.Add($ws.Range("E$row"), $url)
# re-apply the same Hyperlink cell style used by the rest of column E
$ws.Range("E12").Copy()
$ws.Range("E$row").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the current selection to match the author's final cursor position
[void]$ws.Range("D21").Select()
